$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format before writing so numeric-looking
# price strings (e.g. "241.12") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "97.594.73"
$ws.Range("E2").Value = "  +5.70%  "
$ws.Range("D3").Value = "3.120.90"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "241.01"
$ws.Range("E5").Value = "  +2.93%  "
$ws.Range("D6").Value = "611.52"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +3.29%  "
$ws.Range("D8").Value = "0.385"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.115.62"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Value = "0.784"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "97.055.04"
$ws.Range("E13").Value = "  +5.36%  "
$ws.Range("D14").Value = "0.0000241"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "33.93"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "5.37"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "3.708.08"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "3.125.66"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "3.60"
$ws.Range("E19").Value = "  -4.61%  "
$ws.Range("D20").Value = "516.46"
$ws.Range("E20").Value = "  +18.10%  "
$ws.Range("D21").Value = "14.61"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "5.68"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").Value = "0.0000193"
$ws.Range("E23").Value = "  -5.20%  "
$ws.Range("D24").Value = "8.83"
$ws.Range("E24").Value = "  -4.43%  "
$ws.Range("D25").Value = "5.53"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").Value = "86.45"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "11.71"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").Value = "3.290.01"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "0.237"
$ws.Range("E30").Value = "  +3.63%  "
$ws.Range("D31").Value = "0.175"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "0.125"
$ws.Range("E32").Value = "  +6.19%  "
$ws.Range("D33").Value = "9.02"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").Value = "26.72"
$ws.Range("E34").Value = "  +3.82%  "
$ws.Range("D35").Value = "0.846"
$ws.Range("E35").Value = "  -18.48%  "
$ws.Range("D36").Value = "0.151"
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("D37").Value = "7.33"
$ws.Range("E37").Value = "  -8.12%  "
$ws.Range("D38").Value = "495.23"
$ws.Range("E38").Value = "  +6.57%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "24.26"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("B40").Value = "PancakeSwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "0.437"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("D42").Value = "1.24"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").Value = "3.59"
$ws.Range("E43").Value = "  -9.39%  "
$ws.Range("D45").Value = "3.19"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").Value = "162.85"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").Value = "1.93"
$ws.Range("E47").Value = "  +5.30%  "
$ws.Range("D48").Value = "0.695"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "0.0325"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "4.38"
$ws.Range("E51").Value = "  +1.09%  "

# Restore default style on column D (clears the temporary text format)
$ws.Range("D2:D51").Style = "Normal"

